$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

$times = @(
    "2021-10-05 13:39:43.718607",
    "2021-10-05 13:39:43.718619",
    "2021-10-05 13:39:43.718623",
    "2021-10-05 13:39:43.718626",
    "2021-10-05 13:39:43.718629",
    "2021-10-05 13:39:43.718632",
    "2021-10-05 13:39:43.718635",
    "2021-10-05 13:39:43.718637",
    "2021-10-05 13:39:43.718640",
    "2021-10-05 13:39:43.718643",
    "2021-10-05 13:39:43.718646",
    "2021-10-05 13:39:43.718649",
    "2021-10-05 13:39:43.718652",
    "2021-10-05 13:39:43.718654",
    "2021-10-05 13:39:43.718657",
    "2021-10-05 13:39:43.718660",
    "2021-10-05 13:39:43.718663",
    "2021-10-05 13:39:43.718666",
    "2021-10-05 13:39:43.718668",
    "2021-10-05 13:39:43.718671",
    "2021-10-05 13:39:43.718674",
    "2021-10-05 13:39:43.718676",
    "2021-10-05 13:39:43.718679",
    "2021-10-05 13:39:43.718682",
    "2021-10-05 13:39:43.718685",
    "2021-10-05 13:39:43.718688",
    "2021-10-05 13:39:43.718690",
    "2021-10-05 13:39:43.718693",
    "2021-10-05 13:39:43.718696",
    "2021-10-05 13:39:43.718699",
    "2021-10-05 13:39:43.718702",
    "2021-10-05 13:39:43.718704",
    "2021-10-05 13:39:43.718707",
    "2021-10-05 13:39:43.718710",
    "2021-10-05 13:39:43.718715",
    "2021-10-05 13:39:43.718719",
    "2021-10-05 13:39:43.718723",
    "2021-10-05 13:39:43.718728",
    "2021-10-05 13:39:43.718733",
    "2021-10-05 13:39:43.718737",
    "2021-10-05 13:39:43.718742",
    "2021-10-05 13:39:43.718745",
    "2021-10-05 13:39:43.718779"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
